{"js": "// Move-In/Move-Out Inspection Checklist Template update\n// 1) Turn the first paragraph into a Heading 1 with new title text.\n// 2) Replace the info-line paragraphs with templated placeholder fields.\n// 3) Append a 3-column \"Area / Condition / Notes\" checklist table.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst p = paragraphs.items;\n\n// --- Paragraph 0: title -> Heading 1 style, new wording -------------------\np[0].style = \"Heading 1\";\np[0].getRange().insertText(\"Move-In/Move-Out Inspection Checklist\", Word.InsertLocation.replace);\n\n// --- Paragraphs 1-5: swap boilerplate copy for templated fields -----------\np[1].getRange().insertText(\"Tenant Name: {{TENANT_NAME}}\", Word.InsertLocation.replace);\np[2].getRange().insertText(\"Property Address: {{PROPERTY_ADDRESS}}\", Word.InsertLocation.replace);\np[3].getRange().insertText(\"Unit: {{UNIT_NUMBER}}\", Word.InsertLocation.replace);\np[4].getRange().insertText(\"Inspection Type: {{MOVE_IN_OR_OUT}}\", Word.InsertLocation.replace);\np[5].getRange().insertText(\"Inspection Date: {{INSPECTION_DATE}}\", Word.InsertLocation.replace);\n\nawait context.sync();\n\n// --- Append the room-by-room condition table at the end of the body -------\nconst rows = [\n  [\"Area\", \"Condition\", \"Notes\"],\n  [\"Entry / Hallway\", \"{{CONDITION}}\", \"{{NOTES}}\"],\n  [\"Living Room\", \"{{CONDITION}}\", \"{{NOTES}}\"],\n  [\"Kitchen\", \"{{CONDITION}}\", \"{{NOTES}}\"],\n  [\"Bathroom\", \"{{CONDITION}}\", \"{{NOTES}}\"],\n  [\"Bedroom\", \"{{CONDITION}}\", \"{{NOTES}}\"],\n];\n\nfunction cellXml(text) {\n  return `<w:tc><w:tcPr><w:tcW w:type=\"dxa\" w:w=\"2880\"/></w:tcPr><w:p><w:r><w:t>${text}</w:t></w:r></w:p></w:tc>`;\n}\n\nlet rowsXml = \"\";\nfor (const row of rows) {\n  rowsXml += `<w:tr>${row.map(cellXml).join(\"\")}</w:tr>`;\n}\n\nconst tableXml =\n  `<w:tbl><w:tblPr><w:tblW w:type=\"auto\" w:w=\"0\"/><w:tblLook w:firstColumn=\"1\" w:firstRow=\"1\" w:lastColumn=\"0\" w:lastRow=\"0\" w:noHBand=\"0\" w:noVBand=\"1\" w:val=\"04A0\"/></w:tblPr>` +\n  `<w:tblGrid><w:gridCol w:w=\"2880\"/><w:gridCol w:w=\"2880\"/><w:gridCol w:w=\"2880\"/></w:tblGrid>` +\n  rowsXml +\n  `</w:tbl>`;\n\nconst packageXml =\n  `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>` +\n  `<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">` +\n  `<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">` +\n  `<pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">` +\n  `<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>` +\n  `</Relationships></pkg:xmlData></pkg:part>` +\n  `<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">` +\n  `<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>${tableXml}</w:body></w:document></pkg:xmlData>` +\n  `</pkg:part></pkg:package>`;\n\nparagraphs.load(\"items\");\nawait context.sync();\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.getRange(Word.RangeLocation.after).insertOoxml(packageXml, Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Move-In/Move-Out Inspection Checklist Template update\n# 1) Turn the first paragraph into a Heading 1 with new title text.\n# 2) Replace the info-line paragraphs with templated placeholder fields.\n# 3) Append a 3-column \"Area / Condition / Notes\" checklist table.\n\n$d = $word.ActiveDocument\n\n# --- Paragraph 1: title -> Heading 1 style, new wording -------------------\n$p1 = $d.Paragraphs.Item(1)\n$p1.Style = \"Heading 1\"\n$p1.Range.Text = \"Move-In/Move-Out Inspection Checklist\"\n\n# --- Paragraphs 2-6: swap boilerplate copy for templated fields -----------\n$d.Paragraphs.Item(2).Range.Text = \"Tenant Name: {{TENANT_NAME}}\"\n$d.Paragraphs.Item(3).Range.Text = \"Property Address: {{PROPERTY_ADDRESS}}\"\n$d.Paragraphs.Item(4).Range.Text = \"Unit: {{UNIT_NUMBER}}\"\n$d.Paragraphs.Item(5).Range.Text = \"Inspection Type: {{MOVE_IN_OR_OUT}}\"\n$d.Paragraphs.Item(6).Range.Text = \"Inspection Date: {{INSPECTION_DATE}}\"\n\n# --- Append the room-by-room condition table at the end of the body -------\n$rows = @(\n  @(\"Area\", \"Condition\", \"Notes\"),\n  @(\"Entry / Hallway\", \"{{CONDITION}}\", \"{{NOTES}}\"),\n  @(\"Living Room\", \"{{CONDITION}}\", \"{{NOTES}}\"),\n  @(\"Kitchen\", \"{{CONDITION}}\", \"{{NOTES}}\"),\n  @(\"Bathroom\", \"{{CONDITION}}\", \"{{NOTES}}\"),\n  @(\"Bedroom\", \"{{CONDITION}}\", \"{{NOTES}}\")\n)\n\n$rowsXml = \"\"\nforeach ($row in $rows) {\n  $cellsXml = \"\"\n  foreach ($cell in $row) {\n    $cellsXml += \"<w:tc><w:tcPr><w:tcW w:type=`\"dxa`\" w:w=`\"2880`\"/></w:tcPr><w:p><w:r><w:t>$cell</w:t></w:r></w:p></w:tc>\"\n  }\n  $rowsXml += \"<w:tr>$cellsXml</w:tr>\"\n}\n\n$tableXml = \"<w:tbl><w:tblPr><w:tblW w:type=`\"auto`\" w:w=`\"0`\"/><w:tblLook w:firstColumn=`\"1`\" w:firstRow=`\"1`\" w:lastColumn=`\"0`\" w:lastRow=`\"0`\" w:noHBand=`\"0`\" w:noVBand=`\"1`\" w:val=`\"04A0`\"/></w:tblPr><w:tblGrid><w:gridCol w:w=`\"2880`\"/><w:gridCol w:w=`\"2880`\"/><w:gridCol w:w=`\"2880`\"/></w:tblGrid>$rowsXml</w:tbl>\"\n\n$packageXml = \"<?xml version=`\"1.0`\" encoding=`\"UTF-8`\" standalone=`\"yes`\"?><pkg:package xmlns:pkg=`\"http://schemas.microsoft.com/office/2006/xmlPackage`\"><pkg:part pkg:name=`\"/_rels/.rels`\" pkg:contentType=`\"application/vnd.openxmlformats-package.relationships+xml`\"><pkg:xmlData><Relationships xmlns=`\"http://schemas.openxmlformats.org/package/2006/relationships`\"><Relationship Id=`\"rId1`\" Type=`\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument`\" Target=`\"word/document.xml`\"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name=`\"/word/document.xml`\" pkg:contentType=`\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`\"><pkg:xmlData><w:document xmlns:w=`\"http://schemas.openxmlformats.org/wordprocessingml/2006/main`\"><w:body>$tableXml</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\"\n\n$newP = $d.Paragraphs.Add()\n$null = $newP.Range.InsertXML($packageXml)\n"}
